# Article 39 ("Road to success", Surah Al Imran 165-185) goes live:
# fill in row 40 of Sheet1 with its metadata/content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ayats = "Surah Al Imran, 165 - 185"
$tags = "Cycle of our life, Cheering up, Stairs of our life, Meeting the God"
# The source text has an embedded vertical-tab (Excel's "_x000B_" escape)
# between two paragraphs; splice it back in with [char]11.
$content = "h1: Road to success`r`np: We all want some things in our life that are difficult to achieve. We want higher rank, a good home, excelling children and a respectable status in our community. All these things require tremendous amount of work. `r`np: We wake up early in mornings, prepare ourselves for the day, go office, perform tasks to their best. As the time to return home nears, we pack up and move back to our homes. Our family receive us. They always have some news to share. We are served good tea in a cosy place that we tone down our day’s stress. `r`np: All of these cycles are repeated everyday. Waking up to sleeping back into the bed is focused on achieving a place of comfort in our community. We want better future for our children. We work hard, feed them that they stand up to our expectations.`r`np: A perfect life is lived all our lives. We near our retirement to hold on to our children’s shoulders. They put in their best to keep us protected of elderly sickness. They help us walk back again. Sometimes they succeed, other times they fail. The body stops working as it is supposed to work. Eyes do not see properly, hearing lowers, old age hardships encompasses us. `r`np: The death wears us down. We die, every one mourns the loss and they start preparing for their deaths.`r`nh3: The next step on the stair`r`np: We were moving up the ladder of our life. Death meets us and it takes us further in our pursuit. We earned a lot of goodness during our life and we are lowered into our grave. Angels await us, our Creator meets us and gives us glad tidings. We sleep in our grave and wait for the day of resurrection. We are kept in the comfort of our Creator. It is a good part of life again." + [char]11 + "p: We witness the physical reality of an unending world. It calms us down of the worry of falling into another pit. The days where children were failing and it was stressing me, are now over. This time in grave is only meant for us to feel safe in the beds made by our Creator. `r`np: This time keeps on moving forward till we enter the next phase.`r`nh3: Rising up on the day of resurrection`r`np: The judgement day nears. A horn is blown. We wake back up. We are taken to the ground of justice. We are worried. We might fail here and never enter the paradise. Something keeps assuring us of our children’s good deeds. How we managed to feed them halal and make them the leaders of the future. Our children excelled in their lives, brought good name to muslims. Islam prospered due to their contributions. Our communities rose higher for their inputs. `r`np: These relaxing reflections of our past comfort us. We see creator, He is happy with us, we feel better. Our Prophet (PBUH) is here with us and his companions are standing beside him. It is a better world out here. The plane of Arafat feels like a garden from the heaven. The success is kissing us. People who dominated the world of justice are standing in front and around us.`r`np: We are given the reward of our deeds, it is all clear. `r`nh3: Entering the paradise`r`np: The bridge named ‘Sirat’, narrower then a hair awaits us. As we near it, it also broadens up. The crossing is a little easier. Their the heaven awaits us. Khalid (ra), Abu Bakr (ra), Umar (ra), Ali (ra) and their likes are welcoming us. We meet the companions with a pride on our face. Our shoulders broadened, worries dusted, pride filled hearts meet the true warriors of Islam. `r`nh3: We succeed in our life`r`np: The moment of shocks and turbulence is over. We succeed in finding the Paradise. It is where there are rivers flowing under dense green trees. God is above us looking for us. Angels are awaiting orders. Everything which lacked ‘abundance’ back on earth, is now in ‘abundance’. `r`np: This moment will make us forget the failures of earth; the property disputes to our divorces. Life will get back on track and we will be ruled with love of God above us. `r`nh3: Conclusion`r`np: We take small happenings around us too heavy. Our children wrong us and we get angry. Wives are always angry and sometimes we retaliate which leads us to ‘repenting’. All important things are not older than 30 years; the time passed in a moment. `r`np: These tiresome nights sitting on computer won’t worry us on our death. We might meet Allah of of our hearts with all our love.`r`nquote: Every soul will taste death, and you will only be given your [full] compensation on the Day of Resurrection. So he who is drawn away from the Fire and admitted to Paradise has attained [his desire]. And what is the life of this world except the enjoyment of delusion. <br>- Surah Al Imran verse 185`r`np.note: Guys, if you like this project. Please follow this project's page on twitter. <a href=`"https://twitter.com/zakatlists`">Click here to go to the twitter page</a>. Subscribe below to receive updates."

$ws.Range("A40").Value = 39
$ws.Range("B40").Value = 43869
$ws.Range("B40").NumberFormat = "d-mmm-yy"
$ws.Range("C40").Value = $ayats
$ws.Range("F40").Value = $tags
$ws.Range("D40").Value = $content
$ws.Range("E40").Value = "Qasim Ali"

$ws.Rows.Item(40).RowHeight = 409.6

[void]$ws.Range("D40").Select()
